# The commit rotates the four data rows (2-5) of the "Artfynd" sheet:
#   new row 2 <- old row 3
#   new row 3 <- old row 4
#   new row 4 <- old row 5
#   new row 5 <- old row 2
# (the header in row 1 is untouched). We reproduce that by copying whole
# rows through a scratch row far below the used range, clearing each
# destination immediately before pasting into it so that cells which are
# blank in the source don't leave stale leftovers in the destination.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = "A100:AY100"

# Stash old row 2 out of the way before it gets overwritten.
$ws.Range($scratch).Clear()
$ws.Range("A2:AY2").Copy($ws.Range($scratch))

# Shift rows 3, 4, 5 up by one.
$ws.Range("A2:AY2").Clear()
$ws.Range("A3:AY3").Copy($ws.Range("A2:AY2"))

$ws.Range("A3:AY3").Clear()
$ws.Range("A4:AY4").Copy($ws.Range("A3:AY3"))

$ws.Range("A4:AY4").Clear()
$ws.Range("A5:AY5").Copy($ws.Range("A4:AY4"))

# Drop the stashed original row 2 into row 5.
$ws.Range("A5:AY5").Clear()
$ws.Range($scratch).Copy($ws.Range("A5:AY5"))

# Clean up the scratch area.
$ws.Range($scratch).Clear()
